# Daily attendance processing - 2025-12-25 17:30:58
# Normalize the "Recorded By" (column G) value lists: entries that list
# "System"/"system" before a real user are re-ordered so the real
# user identity comes first, with "System" (or the stray lowercase
# "system") moved after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-value lookup table: old combined "Recorded By" string -> new one.
$map = @{
    "system, backup@backdoor.com, System" = "backup@backdoor.com, System, system";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
